$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (Price column D, Volume(1h)
# column E) to the latest scrape, and swap the WhiteBITCoin/EnergySwap rows
# (50/51) to their new ranking order with updated figures.
#
# Price-column values that look like plain decimal numbers (e.g. "544.16")
# are written with a leading apostrophe, Excel's standard "force text" quote
# prefix, so they stay text cells (matching the sheet's existing inlineStr
# text cells) instead of being auto-converted to numeric values.

$ws.Range('D2').Value = '58.581.92'
$ws.Range('E2').Value = '  -2.39%  '
$ws.Range('D3').Value = '2.290.70'
$ws.Range('E3').Value = '  -5.06%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''544.16'
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('D6').Value = '''130.81'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''0.571'
$ws.Range('E8').Value = '  -2.95%  '
$ws.Range('D9').Value = '2.289.17'
$ws.Range('E9').Value = '  -5.11%  '
$ws.Range('E10').Value = '  -3.33%  '
$ws.Range('D11').Value = '''5.50'
$ws.Range('E11').Value = '  -2.85%  '
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('E13').Value = '  -5.76%  '
$ws.Range('D14').Value = '''23.89'
$ws.Range('E14').Value = '  -3.36%  '
$ws.Range('D15').Value = '2.701.55'
$ws.Range('E15').Value = '  -5.04%  '
$ws.Range('D16').Value = '58.566.30'
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('D17').Value = '''0.0000132'
$ws.Range('E17').Value = '  -3.75%  '
$ws.Range('D18').Value = '2.253.26'
$ws.Range('E18').Value = '  -8.03%  '
$ws.Range('D19').Value = '''10.60'
$ws.Range('E19').Value = '  -5.36%  '
$ws.Range('D20').Value = '''4.30'
$ws.Range('E20').Value = '  -5.15%  '
$ws.Range('D21').Value = '''313.60'
$ws.Range('E22').Value = '  -5.36%  '
$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = '''63.28'
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('D25').Value = '''0.167'
$ws.Range('E25').Value = '  -6.67%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = '''8.08'
$ws.Range('E27').Value = '  -6.41%  '
$ws.Range('D28').Value = '''1.32'
$ws.Range('E28').Value = '  -5.86%  '
$ws.Range('D29').Value = '''1.74'
$ws.Range('E29').Value = '  -1.99%  '
$ws.Range('D30').Value = '''170.41'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').Value = '0.0₃0723'
$ws.Range('E31').Value = '  -6.14%  '
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('D33').Value = '''5.76'
$ws.Range('E33').Value = '  -5.70%  '
$ws.Range('D34').Value = '''0.376'
$ws.Range('E34').Value = '  -6.07%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  -3.75%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '''1.24'
$ws.Range('E38').Value = '  -6.00%  '
$ws.Range('D39').Value = '''3.95'
$ws.Range('E39').Value = '  -6.22%  '
$ws.Range('D40').Value = '''38.02'
$ws.Range('E40').Value = '  -2.28%  '
$ws.Range('E41').Value = '  -5.67%  '
$ws.Range('D42').Value = '''297.93'
$ws.Range('E42').Value = '  -8.34%  '
$ws.Range('D43').Value = '''139.95'
$ws.Range('E43').Value = '  -4.13%  '
$ws.Range('E44').Value = '  -5.03%  '
$ws.Range('D45').Value = '''0.0949'
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('D46').Value = '''0.0497'
$ws.Range('E46').Value = '  -3.41%  '
$ws.Range('D47').Value = '''0.552'
$ws.Range('E47').Value = '  -4.35%  '
$ws.Range('D48').Value = '''18.41'
$ws.Range('E48').Value = '  -6.96%  '
$ws.Range('E49').Value = '  -3.46%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').Value = '''11.02'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''16.52'
$ws.Range('E51').Value = '  -4.78%  '
